$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Resize / reposition the workbook window (mirrors workbookView xWindow/yWindow/windowWidth)
$excel.ActiveWindow.Left = 120
$excel.ActiveWindow.Top = 0
$excel.ActiveWindow.Width = 12540

# Replace STATUS "TODO" with "DONE" everywhere it currently appears. Using
# Replace() mutates the single shared string in place instead of minting a
# brand-new string entry, matching how Excel's Find & Replace collapses
# every occurrence of the old text onto the same shared-string slot.
$ws.Range("E2:E12").Replace("TODO", "DONE")

# New backlog items (rows 13-16). Order matters: it controls the order new
# unique strings are appended to the shared-string table.
$ws.Range("A13").Value = "PHP background"

# Mark a handful of existing rows as REMOVED.
$ws.Range("E5").Value = "REMOVED"

$ws.Range("A14").Value = "PHP tekst"
$ws.Range("A15").Value = "PHP MVC system"
$ws.Range("A16").Value = "PHP functions"

$ws.Range("E6").Value = "REMOVED"
$ws.Range("E8").Value = "REMOVED"
$ws.Range("E9").Value = "REMOVED"
$ws.Range("E10").Value = "REMOVED"

# Fill in the rest of the new backlog rows (13-16) -- priority / size / fase / status
$ws.Range("B13").Value = "HIGH"
$ws.Range("C13").Value = "XL"
$ws.Range("D13").Value = 2
$ws.Range("E13").Value = "DONE"

$ws.Range("B14").Value = "HIGH"
$ws.Range("C14").Value = "XL"
$ws.Range("D14").Value = 2
$ws.Range("E14").Value = "DONE"

$ws.Range("B15").Value = "MAX"
$ws.Range("C15").Value = "XL"
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = "DONE"

$ws.Range("B16").Value = "HIGH"
$ws.Range("C16").Value = "XL"
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = "DONE"

# Move the active cell selection to match the saved cursor position
$ws.Range("B18").Select()
